$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.548.47"
$ws.Range("E2").Value = "  +1.42%  "

# Row 3
$ws.Range("D3").Value = "3.357.29"
$ws.Range("E3").Value = "  +1.12%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.80%  "

# Row 8
$ws.Range("E8").Value = "  +0.70%  "

# Row 9
$ws.Range("E9").Value = "  +4.47%  "

# Row 10
$ws.Range("E10").Value = "  +1.37%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "47.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.93%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "693.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.65%  "

# Row 14
$ws.Range("D14").Value = "3.911.99"
$ws.Range("E14").Value = "  +1.15%  "

# Row 15
$ws.Range("E15").Value = "  +1.30%  "

# Row 16
$ws.Range("D16").Value = "68.622.25"
$ws.Range("E16").Value = "  +1.48%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.392.21"
$ws.Range("E17").Value = "  +1.98%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.120"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.55%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.97%  "

# Row 21
$ws.Range("E21").Value = "  +1.38%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.43%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.63%  "

# Row 25
$ws.Range("E25").Value = "  +2.08%  "

# Row 26
$ws.Range("E26").Value = "  +1.91%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.33%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.87%  "

# Row 29
$ws.Range("E29").Value = "  +1.90%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.25%  "

# Row 31
$ws.Range("E31").Value = "  +1.87%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "547.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.51%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.08%  "

# Row 36
$ws.Range("D36").Value = "3.699.89"
$ws.Range("E36").Value = "  +0.90%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.50%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.142"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.18%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.45%  "

# Row 40
$ws.Range("E40").Value = "  +2.86%  "

# Row 41
$ws.Range("E41").Value = "  +0.38%  "

# Row 42
$ws.Range("D42").Value = "0.0₃0674"
$ws.Range("E42").Value = "  +2.12%  "

# Row 43
$ws.Range("E43").Value = "  +1.14%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0414"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.20%  "

# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.86%  "

# Row 46
$ws.Range("E46").Value = "  +2.00%  "

# Row 47
$ws.Range("E47").Value = "  +0.98%  "

# Row 48
$ws.Range("E48").Value = "  -0.25%  "

# Row 49
$ws.Range("E49").Value = "  -1.38%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.67%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.56%  "
